$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price records (rows 2-5) get their Fecha/Volumen/Precio
# columns rotated down one row (row2 <- row5, row3 <- row4, row4 <- row2,
# row5 <- row3), while the rest of each row's data stays put.

$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16500
$ws.Range("S2").Value = 825

$ws.Range("D3").Value = 44533
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("S3").Value = 825

$ws.Range("D4").Value = 44708
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("S4").Value = 1025

$ws.Range("D5").Value = 44357
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 725
